$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "apple" entry was renamed to "test" (A2), which is the only real
# content edit in the diff -- the shifted shared-string indices elsewhere
# in the sheet are just a side effect of "apple" dropping out of the
# shared-strings table and "test" being appended at the end.
$ws.Range("A2").Value = "test"

# Reflect the new active selection recorded in the sheet view.
$ws.Range("A2").Select()
